$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New values for columns B, C, D, E, G across rows 2-11
$data = @{
    2  = @{ B = 0.6606524410359556;   C = 1.655778082260271;   D = 3.537761648806719;   E = 0.4942365360607697;   G = 6.348428708163715 }
    3  = @{ B = 3.286832544864788;    C = 1.655778082260271;   D = 0.1494219747398047;  E = 0.4942365360607697;   G = 5.586269137925634 }
    4  = @{ B = 3.286832544864788;    C = 1.655778082260271;   D = 0.1494219747398047;  E = 10.19245300693656;    G = 15.28448560880142 }
    5  = @{ B = 0.2917716402565462;   C = 1.655778082260271;   D = 0.1494219747398047;  E = 0.4942365360607697;   G = 2.591208233317391 }
    6  = @{ B = 0.000001295275857016165; C = 0.04071648406533734; D = 0.1494219747398047;  E = 0.4942365360607697;   G = 0.6843762901417687 }
    7  = @{ B = 3.286832544864788;    C = 1.655778082260271;   D = 0.7527432677738641;  E = 0.4942365360607697;   G = 6.189590430959694 }
    8  = @{ B = 1.455362044514542;    C = 1.655778082260271;   D = 22.3905356188092;    E = 10.19245300693656;    G = 35.69412875252057 }
    9  = @{ B = 0.04271373187048222;  C = 0.04071648406533734; D = 0.7527432677738641;  E = 0.4942365360607697;   G = 1.330410019770453 }
    10 = @{ B = 1.455362044514542;    C = 1.655778082260271;   D = 3.537761648806719;   E = 0.4942365360607697;   G = 7.143138311642302 }
    11 = @{ B = 3.286832544864788;    C = 1.655778082260271;   D = 0.7527432677738641;  E = 0.4942365360607697;   G = 6.189590430959694 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
